# Apply the "Updated WBS" edit to the Traffic Sign Recognition capstone
# workbook: the WBS worksheet is refilled with the new project's task
# breakdown (Dark_Web_Crusaders / Traffic Sign Recognition System),
# replacing the generic placeholder template rows, and the two
# illustrative callout shapes are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WBS")

function Set-Text($addr, $text) {
    $ws.Range($addr).Value2 = $text
}

function Set-Day($addr, $year, $month, $day) {
    $ws.Range($addr).Value2 = (Get-Date -Year $year -Month $month -Day $day -Hour 0 -Minute 0 -Second 0)
}

# --- Header area ---------------------------------------------------------
Set-Text "J4" "In-progress"
Set-Text "E8" "Dark_Web_Crusaders"
Set-Text "E9" "Traffic Sign Recongition System Using CNN"

# --- Row 14: Phase 1 ------------------------------------------------------
Set-Text "B14" "1. Exploring Generalized Traffic Sign Dataset "
Set-Text "F14" "28/9/2022"
Set-Text "G14" "28/9/2022"
Set-Text "H14" "Shamaem "

# --- Row 15: 2. Data Acquisition and EDA ----------------------------------
Set-Text "B15" "2. Data Acquisition and EDA"
$ws.Range("C15").ClearContents()
Set-Day "F15" 2022 1 10
Set-Day "G15" 2022 3 10
Set-Text "H15" "Shamaem "

# --- Row 16: 2.1 Data Acquisition -----------------------------------------
Set-Text "C16" "2.1  Data Acquisition"
$ws.Range("D16").ClearContents()
Set-Day "F16" 2022 1 10
Set-Day "G16" 2022 1 10
Set-Text "H16" "Shamaem "

# --- Row 17: 2.2 EDA --------------------------------------------------------
Set-Text "C17" "2.2  EDA"
$ws.Range("E17").ClearContents()
Set-Day "F17" 2022 2 10
Set-Day "G17" 2022 3 10
Set-Text "H17" "Shamaem "

# --- Row 18: 3 Pre-Processing ----------------------------------------------
Set-Text "B18" "3  Pre-Processing"
Set-Day "F18" 2022 4 10
Set-Day "G18" 2022 5 10
Set-Text "H18" "Shamaem "
Set-Text "J18" "Done"

# --- Row 19: 3,1 Normalization ----------------------------------------------
Set-Text "C19" "3,1   Normalization"
Set-Day "F19" 2022 4 10
Set-Day "G19" 2022 4 10
Set-Text "H19" "Shamaem "
Set-Text "J19" "Done"

# --- Row 20: 3.2 Data Augmentation -------------------------------------------
Set-Text "C20" "3.2  Data Augmentation"
Set-Day "F20" 2022 4 10
Set-Day "G20" 2022 4 10
Set-Text "H20" "Shamaem "
Set-Text "J20" "Done"

# --- Row 21: 3.3 EDA after pre-processing ------------------------------------
Set-Text "C21" "3.3 EDA after pre-processing"
Set-Day "F21" 2022 5 10
Set-Day "G21" 2022 5 10
Set-Text "H21" "Shamaem "
Set-Text "J21" "Done"

# --- Row 22: 4 Model -----------------------------------------------------------
Set-Text "B22" "4  Model"
Set-Day "F22" 2022 6 10
Set-Day "G22" 2022 10 10
Set-Text "H22" "Toheed"
Set-Text "J22" "Done"

# --- Row 23: 4.1 Building The Model ---------------------------------------------
Set-Text "C23" "4.1 Building The Model"
Set-Day "F23" 2022 6 10
Set-Day "G23" 2022 6 10
Set-Text "H23" "Toheed"
Set-Text "J23" "Done"

# --- Row 24: 4.2 Training The Model ----------------------------------------------
Set-Text "C24" "4.2 Training The Model"
Set-Day "F24" 2022 7 10
Set-Day "G24" 2022 7 10
Set-Text "H24" "Toheed"
Set-Text "J24" "Done"

# --- Row 25: 4.3 Validating The Model --------------------------------------------
Set-Text "C25" "4.3 Validating The Model"
Set-Day "F25" 2022 8 10
Set-Day "G25" 2022 8 10
Set-Text "H25" "Toheed"
Set-Text "J25" "Done"

# --- Row 26: 4.4 Testing and Evulation -------------------------------------------
Set-Text "C26" "4.4  Testing and Evulation"
Set-Day "F26" 2022 9 10
Set-Day "G26" 2022 9 10
Set-Text "H26" "Toheed"
Set-Text "J26" "Done"

# --- Row 27: 4.5 Visulaizaton and Results ----------------------------------------
Set-Text "C27" "4.5 Visulaizaton and Results"
Set-Day "F27" 2022 10 10
Set-Day "G27" 2022 10 10
Set-Text "H27" "Toheed"
Set-Text "J27" "Done"

# --- Row 28: 5 Building Interface -------------------------------------------------
Set-Text "B28" "5  Building Interface"
Set-Day "F28" 2022 11 10
Set-Day "G28" 2022 12 10
Set-Text "H28" "Shamaem /Toheed"
Set-Text "J28" "Delayed"

# --- Row 29: clear the now-unused status ------------------------------------------
$ws.Range("J29").ClearContents()

# --- Remove the two illustrative callout shapes -----------------------------------
$shapeCount = $ws.Shapes.Count
for ($i = $shapeCount; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# --- Update the remembered selection on the WBS sheet -----------------------------
$ws.Range("I25").Select()
